$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto data.
# D-column values are forced to Text format before assignment so that Excel
# does not reinterpret numeric-looking strings (e.g. "1.000", "0.02410") as
# numbers and strip significant trailing zeros / re-format them. The cell
# style is reset back to "Normal" afterwards so no stray style index is left
# attached to the cell (keeping it identical in shape to the original file).

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '29.966.33'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.49%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.866.70'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -2.99%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '317.93'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -2.68%  '

$ws.Range("E6").Value = '  +0.09%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5084'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -1.47%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.3911'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -2.41%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.08155'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -3.65%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '41.92'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -2.31%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '1.088'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -3.08%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '22.63'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +6.46%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.875.21'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -2.33%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '6.252'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -1.35%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '7.143'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -2.72%  '

$ws.Range("E16").Value = '  +0.17%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '91.72'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -2.75%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '0.00001075'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -3.74%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.06341'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -6.43%  '

$ws.Range("E20").Value = '  -1.32%  '

$ws.Range("E21").Value = '  +0.08%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '29.943.49'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.56%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '5.783'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -4.60%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '11.06'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -1.30%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.202'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.084.98'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -2.54%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '160.29'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.04%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '20.83'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.72%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '2.219'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -10.15%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '126.30'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -1.98%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.1032'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.55%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.039'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -3.48%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.843'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.78%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '3.734'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +2.21%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.02410'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -3.68%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.06318'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -4.34%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '5.169'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.66%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.2137'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -4.02%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.168'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -6.06%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '8.449'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -6.24%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.6253'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -4.46%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.205'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -3.31%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '11.22'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -1.45%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.04%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.5860'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -4.55%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '12.83'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -2.83%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '3.622'
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.984'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -3.50%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '121.89'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -3.08%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.199'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -3.57%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.150'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.15%  '
